# Apply the commit: "PageUtility SubcategoryClass ManageUserClass EdittestcaseinCategoryClass"
#  - Login sheet: add two new rows of test credentials
#  - NewsSearch sheet: rename one test data entry
#  - Add new "Subcategory" sheet with test data
#  - Category sheet: rename one test data entry
#  - Add new "ManageUsers" sheet with test data
#  - Update selections / active tab to match the new recorded Selenium state

$wb = $excel.ActiveWorkbook

$wsLogin = $wb.Worksheets.Item("Login")
$wsCategory = $wb.Worksheets.Item("Category")
$wsNews = $wb.Worksheets.Item("NewsSearch")

# --- Sheet "Login": append rows 6 and 7 with new admin credentials ---
$wsLogin.Range("A6").Value = "admin123"
$wsLogin.Range("B6").Value = "test"
$wsLogin.Range("A7").Value = "admin435"
$wsLogin.Range("B7").Value = "test456"
$wsLogin.Range("C10").Select()

# --- Sheet "NewsSearch": row 3 test value changed from "R" to "news" ---
$wsNews.Range("A3").Value = "news"
$wsNews.Range("C4").Select()

# --- New sheet "Subcategory" added after the last existing sheet (NewsSearch) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsSubcategory = $wb.Worksheets.Add($null, $lastSheet)
$wsSubcategory.Name = "Subcategory"
$wsSubcategory.Range("A1").Value = "SubcategoryName"
$wsSubcategory.Range("A2").Value = "AddToyForTest_2190876"
# target column width is 22.109375 characters; engine quantizes to 1/6, use nearest achievable input
$wsSubcategory.Columns.Item(1).ColumnWidth = 21.3
$wsSubcategory.Range("J20").Select()

# --- Sheet "Category": row 3 test value changed to "CategoryEdited452" ---
$wsCategory.Range("A3").Value = "CategoryEdited452"
# target column width is 17.77734375 characters; engine quantizes to 1/6, use nearest achievable input
$wsCategory.Columns.Item(1).ColumnWidth = 16.95

# --- New sheet "ManageUsers" added after "Subcategory" ---
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsManageUsers = $wb.Worksheets.Add($null, $lastSheet2)
$wsManageUsers.Name = "ManageUsers"
$wsManageUsers.Range("A1").Value = "Username"
$wsManageUsers.Range("B1").Value = "Password"
$wsManageUsers.Range("B2").Value = "testfordemo503"
$wsManageUsers.Range("A2").Value = "saumya1993"
# target column widths are 12.21875 / 14.44140625 characters; engine quantizes to 1/6, use nearest achievable input
$wsManageUsers.Columns.Item(1).ColumnWidth = 11.3
$wsManageUsers.Columns.Item(2).ColumnWidth = 13.65
$wsManageUsers.Range("C8").Select()

# --- Final active tab should be "Category" ---
$wsCategory.Range("A3").Select()
